$wb = $excel.ActiveWorkbook

# Sheet "Overview": row 3 corresponds to d3f61a6d-25c7-48a1-bf99-9242a1d9b79d.md
# Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

# Sheet "zh-cn": row 3 is the d3f61a6d file - status + handback datetime updated
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("H3").Value = "2016-03-11 09:36:29"

# Sheet "de-de": row 3 is the d3f61a6d file - status + handback datetime updated
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("H3").Value = "2016-03-11 09:36:34"
